# Update LR-pair worksheet (Mmrn2-Cd248) with values recomputed from new TPM input.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Resolving-Mac" sending-cluster block (previously rows 14-17) is absent
# from the regenerated output, so the table shrinks from 17 to 13 data+header rows.
$ws.Rows("14:17").Delete()

# Recalculated metric columns (E:T) for the remaining rows 2-13.
# Columns A:D (Sending cluster / Ligand symbol / Receptor symbol / Target cluster)
# are unchanged.
$newValues = @{
    2 = @(3, 1, 53.20407366666667, 159.612221, 0.9416818042520982, 0.9416818042520981, 3, 1, 0.7661790000000001, 2.298537, 0.01101161754618133, 0.01101161754618133, 40.76384395785301, 366.874595620677, 0.0103694398786221, 0.01036943987862209)
    3 = @(3, 1, 53.20407366666667, 159.612221, 0.9416818042520982, 0.9416818042520981, 3, 1, 59.42247633333333, 178.267429, 0.8540270394120407, 0.8540270394120407, 3161.517808294423, 28453.66027464981, 0.8042217233536083, 0.8042217233536082)
    4 = @(3, 1, 53.20407366666667, 159.612221, 0.9416818042520982, 0.9416818042520981, 3, 1, 9.295765000000001, 27.887295, 0.1335998624070593, 0.1335998624070593, 494.5725658480217, 4451.153092632196, 0.1258085594793117, 0.1258085594793117)
    5 = @(3, 1, 53.20407366666667, 159.612221, 0.9416818042520982, 0.9416818042520981, 2, 0.6666666666666666, 0.09473066666666667, 0.284192, 0.001361480634718677, 0.001361480634718677, 5.040057367825778, 45.360516310432, 0.001282081540556176, 0.001282081540556176)
    6 = @(3, 1, 1.641018666666667, 4.923056, 0.02904509584209167, 0.02904509584209167, 3, 1, 0.7661790000000001, 2.298537, 0.01101161754618133, 0.01101161754618133, 1.257314041008, 11.315826369072, 0.000319833487005295, 0.000319833487005295)
    7 = @(3, 1, 1.641018666666667, 4.923056, 0.02904509584209167, 0.02904509584209167, 3, 1, 59.42247633333333, 178.267429, 0.8540270394120407, 0.8540270394120407, 97.51339288255821, 877.6205359430239, 0.02480529721146053, 0.02480529721146053)
    8 = @(3, 1, 1.641018666666667, 4.923056, 0.02904509584209167, 0.02904509584209167, 3, 1, 9.295765000000001, 27.887295, 0.1335998624070593, 0.1335998624070593, 15.25452388594667, 137.29071497352, 0.003880420808103298, 0.003880420808103298)
    9 = @(3, 1, 1.641018666666667, 4.923056, 0.02904509584209167, 0.02904509584209167, 2, 0.6666666666666666, 0.09473066666666667, 0.284192, 0.001361480634718677, 0.001361480634718677, 0.1554547923057778, 1.399093130752, 0.00003954433552255578, 0.00003954433552255578)
    10 = @(3, 1, 1.653900666666667, 4.961702, 0.02927309990581012, 0.02927309990581012, 3, 1, 0.7661790000000001, 2.298537, 0.01101161754618133, 0.01101161754618133, 1.267183958886, 11.404655629974, 0.0003223441805539378, 0.0003223441805539376)
    11 = @(3, 1, 1.653900666666667, 4.961702, 0.02927309990581012, 0.02927309990581012, 3, 1, 59.42247633333333, 178.267429, 0.8540270394120407, 0.8540270394120407, 98.27887322268421, 884.509859004158, 0.0250000188469719, 0.0250000188469719)
    12 = @(3, 1, 1.653900666666667, 4.961702, 0.02927309990581012, 0.02927309990581012, 3, 1, 9.295765000000001, 27.887295, 0.1335998624070593, 0.1335998624070593, 15.37427193067667, 138.36844737609, 0.003910882119644333, 0.003910882119644332)
    13 = @(3, 1, 1.653900666666667, 4.961702, 0.02927309990581012, 0.02927309990581012, 2, 0.6666666666666666, 0.09473066666666667, 0.284192, 0.001361480634718677, 0.001361480634718677, 0.1566751127537778, 1.410076014784, 0.00003985475863994561, 0.0000398547586399456)
}

foreach ($rowNum in $newValues.Keys) {
    $vals = $newValues[$rowNum]
    $col = 5
    foreach ($v in $vals) {
        $ws.Cells.Item($rowNum, $col).Value = $v
        $col = $col + 1
    }
}

